$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted as the first data row of this block
# (row 60), pushing the existing rows 60:137 down to 61:138.
$ws.Rows("60:60").Insert()

# Populate the new row 60 with the new record's data.
$ws.Range("A60").Value = 5
$ws.Range("B60").Value = 'Macroferia Regional de Talca'
$ws.Range("C60").Value = 'Maule'
$ws.Range("D60").Value = 44763
$ws.Range("E60").Value = 7
$ws.Range("F60").Value = 'Fruta'
$ws.Range("G60").Value = 100108
$ws.Range("H60").Value = 'Tropicales y subtropicales'
$ws.Range("I60").Value = 100108002
$ws.Range("J60").Value = 'Mango'
$ws.Range("K60").Value = 'Sin especificar'
$ws.Range("L60").Value = 'Primera'
$ws.Range("M60").Value = 228
$ws.Range("N60").Value = 8000
$ws.Range("O60").Value = 8000
$ws.Range("P60").Value = 8000
$ws.Range("Q60").Value = '$/bandeja 4 kilos'
$ws.Range("R60").Value = 'Brasil'
$ws.Range("S60").Value = 2000
$ws.Range("T60").Value = 4
